$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-07 Wednesday" "2024-08-08 Thursday"

Replace-Text "542÷7=77, 3" "810÷2=405, 0"
Replace-Text "726÷3=242, 0" "748÷4=187, 0"
Replace-Text "854÷9=94, 8" "145÷9=16, 1"
Replace-Text "696÷5=139, 1" "831÷7=118, 5"
Replace-Text "951÷7=135, 6" "417÷5=83, 2"

Replace-Text "637÷6=106, 1" "252÷7=36, 0"
Replace-Text "788÷8=98, 4" "352÷5=70, 2"
Replace-Text "422÷3=140, 2" "787÷3=262, 1"
Replace-Text "809÷7=115, 4" "428÷7=61, 1"
Replace-Text "848÷7=121, 1" "615÷2=307, 1"

Replace-Text "728÷2=364, 0" "991÷6=165, 1"
Replace-Text "840÷3=280, 0" "568÷4=142, 0"
Replace-Text "946÷9=105, 1" "663÷2=331, 1"
Replace-Text "505÷6=84, 1" "830÷9=92, 2"
Replace-Text "141÷6=23, 3" "955÷3=318, 1"

Replace-Text "931÷6=155, 1" "533÷7=76, 1"
Replace-Text "708÷3=236, 0" "655÷4=163, 3"
Replace-Text "803÷5=160, 3" "459÷6=76, 3"
Replace-Text "163÷3=54, 1" "883÷7=126, 1"
Replace-Text "773÷9=85, 8" "563÷7=80, 3"

Replace-Text "724÷9=80, 4" "837÷6=139, 3"
Replace-Text "689÷3=229, 2" "120÷5=24, 0"
Replace-Text "815÷5=163, 0" "733÷5=146, 3"
Replace-Text "977÷6=162, 5" "868÷5=173, 3"
Replace-Text "132÷3=44, 0" "561÷5=112, 1"
